# Updates the table style used by the three summary tables (slides 14-16)
# from the custom "Table_0" style to the built-in table style
# {F804C0DF-5358-4515-B3EA-FB9316C64E97}.

$p = $ppt.ActivePresentation

$oldStyleId = "{3044271A-F251-4437-835B-A7D3AC658DE0}"
$newStyleId = "{F804C0DF-5358-4515-B3EA-FB9316C64E97}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
